# =====================================================================
# Import_Template sheet: widen the dataset from 6 example rows to 18,
# restyle the header row, freeze it, and retune the column widths.
# =====================================================================

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Import_Template")

# ---- Column widths ----------------------------------------------------
# ColumnWidth is in "characters"; the runtime round-trips it through a
# pixel conversion that adds back 5/6 of a character, so each desired
# final width is entered short by that amount to land on the integer we
# actually want stored in the sheet.
$colWidths = @(28, 16, 12, 42, 12, 10, 10, 26, 18, 34, 10)
for ($c = 1; $c -le $colWidths.Length; $c++) {
  $ws.Columns.Item($c).ColumnWidth = $colWidths[$c - 1] - (5.0/6.0)
}

# ---- Header row style ---------------------------------------------------
# The bold/dark header fill (style index 1) is removed; header cells fall
# back to the workbook default "Normal" style.
$ws.Range("A1:K1").Style = "Normal"

# ---- Freeze header row ----------------------------------------------------
$ws.Range("A2").Select() | Out-Null
($excel.ActiveWindow.FreezePanes = $true) | Out-Null

# ---- Data rows (2-19) ------------------------------------------------
# Helper: write a single cell, forcing a leading apostrophe for values that
# look like dates (DOB column) so Excel keeps them as literal text instead
# of converting to a date serial.

# Row 2
$ws.Cells.Item(2, 1).Value = 'Andy JONES'
$ws.Cells.Item(2, 2).Value = 'Person'
$ws.Cells.Item(2, 3).Value = "'" + '01/01/2000'
$ws.Cells.Item(2, 4).Value = '1 Osmond Drive, Wells, Somerset'
$ws.Cells.Item(2, 5).Value = 'BA5 2JX'
$ws.Cells.Item(2, 6).ClearContents()
$ws.Cells.Item(2, 7).ClearContents()
$ws.Cells.Item(2, 8).Value = 'ACME Logistics Ltd'
$ws.Cells.Item(2, 9).Value = 'Associate'
$ws.Cells.Item(2, 10).Value = 'Primary subject'
$ws.Cells.Item(2, 11).Value = 'ET5'

# Row 3
$ws.Cells.Item(3, 1).Value = 'Sarah PATEL'
$ws.Cells.Item(3, 2).Value = 'Person'
$ws.Cells.Item(3, 3).Value = "'" + '14/08/1987'
$ws.Cells.Item(3, 4).Value = 'Flat 2, 18 Broad Street, Bristol'
$ws.Cells.Item(3, 5).Value = 'BS1 2HG'
$ws.Cells.Item(3, 6).ClearContents()
$ws.Cells.Item(3, 7).ClearContents()
$ws.Cells.Item(3, 8).Value = 'ACME Logistics Ltd'
$ws.Cells.Item(3, 9).Value = 'Director/Officer'
$ws.Cells.Item(3, 10).Value = 'Officer linkage'
$ws.Cells.Item(3, 11).Value = 'ET5'

# Row 4
$ws.Cells.Item(4, 1).Value = 'Danny KENT'
$ws.Cells.Item(4, 2).Value = 'Person'
$ws.Cells.Item(4, 3).Value = "'" + '23/03/1983'
$ws.Cells.Item(4, 4).Value = 'Terminal 3, Manchester Airport'
$ws.Cells.Item(4, 5).Value = 'M90 1QX'
$ws.Cells.Item(4, 6).ClearContents()
$ws.Cells.Item(4, 7).ClearContents()
$ws.Cells.Item(4, 8).Value = 'KLM1142'
$ws.Cells.Item(4, 9).Value = 'Passenger'
$ws.Cells.Item(4, 10).Value = 'Travel association'
$ws.Cells.Item(4, 11).Value = 'ET5'

# Row 5
$ws.Cells.Item(5, 1).Value = 'Paul SHARP'
$ws.Cells.Item(5, 2).Value = 'Person'
$ws.Cells.Item(5, 3).Value = "'" + '11/11/1979'
$ws.Cells.Item(5, 4).Value = 'Regent''s Park, London'
$ws.Cells.Item(5, 5).Value = 'NW1 4NR'
$ws.Cells.Item(5, 6).ClearContents()
$ws.Cells.Item(5, 7).ClearContents()
$ws.Cells.Item(5, 8).Value = 'North Dock Warehouse'
$ws.Cells.Item(5, 9).Value = 'Visited'
$ws.Cells.Item(5, 10).Value = 'Seen at location'
$ws.Cells.Item(5, 11).Value = 'ET5'

# Row 6
$ws.Cells.Item(6, 1).Value = 'ACME Logistics Ltd'
$ws.Cells.Item(6, 2).Value = 'Organisation'
$ws.Cells.Item(6, 3).ClearContents()
$ws.Cells.Item(6, 4).Value = '1 Canada Square, London'
$ws.Cells.Item(6, 5).Value = 'E14 5AB'
$ws.Cells.Item(6, 6).ClearContents()
$ws.Cells.Item(6, 7).ClearContents()
$ws.Cells.Item(6, 8).Value = 'North Dock Warehouse'
$ws.Cells.Item(6, 9).Value = 'Owns/Uses'
$ws.Cells.Item(6, 10).Value = 'UK logistics business'
$ws.Cells.Item(6, 11).Value = 'ET4'

# Row 7
$ws.Cells.Item(7, 1).Value = 'Moonshadow Marine Ltd'
$ws.Cells.Item(7, 2).Value = 'Organisation'
$ws.Cells.Item(7, 3).ClearContents()
$ws.Cells.Item(7, 4).Value = 'Quayside House, Liverpool'
$ws.Cells.Item(7, 5).Value = 'L3 1BP'
$ws.Cells.Item(7, 6).ClearContents()
$ws.Cells.Item(7, 7).ClearContents()
$ws.Cells.Item(7, 8).Value = 'Moonshadow'
$ws.Cells.Item(7, 9).Value = 'Operates'
$ws.Cells.Item(7, 10).Value = 'Maritime operator'
$ws.Cells.Item(7, 11).Value = 'ET4'

# Row 8
$ws.Cells.Item(8, 1).Value = 'North Dock Warehouse'
$ws.Cells.Item(8, 2).Value = 'Location'
$ws.Cells.Item(8, 4).Value = 'North Dock Road, London'
$ws.Cells.Item(8, 5).Value = 'E16 2GT'
$ws.Cells.Item(8, 8).Value = 'Andy JONES'
$ws.Cells.Item(8, 9).Value = 'Meeting Location'
$ws.Cells.Item(8, 10).Value = 'Regular meetup site'
$ws.Cells.Item(8, 11).Value = 'ET1'

# Row 9
$ws.Cells.Item(9, 1).Value = '500 Silver Street, Greensville'
$ws.Cells.Item(9, 2).Value = 'Location'
$ws.Cells.Item(9, 4).Value = '500 Silver Street, Greensville'
$ws.Cells.Item(9, 5).Value = 'SW1A 1AA'
$ws.Cells.Item(9, 8).Value = 'Sarah PATEL'
$ws.Cells.Item(9, 9).Value = 'Residence'
$ws.Cells.Item(9, 10).Value = 'Address from chart'
$ws.Cells.Item(9, 11).Value = 'ET1'

# Row 10
$ws.Cells.Item(10, 1).Value = 'KLM1142'
$ws.Cells.Item(10, 2).Value = 'Communication'
$ws.Cells.Item(10, 6).Value = 51.47
$ws.Cells.Item(10, 7).Value = -0.4543
$ws.Cells.Item(10, 8).Value = 'Manchester Airport (MAN)'
$ws.Cells.Item(10, 9).Value = 'Route'
$ws.Cells.Item(10, 10).Value = 'Flight identifier'
$ws.Cells.Item(10, 11).Value = 'ET8'

# Row 11
$ws.Cells.Item(11, 1).Value = 'BAW23T'
$ws.Cells.Item(11, 2).Value = 'Communication'
$ws.Cells.Item(11, 6).Value = 50.939
$ws.Cells.Item(11, 7).Value = -1.404
$ws.Cells.Item(11, 8).Value = 'CULDROSE (EGDR)'
$ws.Cells.Item(11, 9).Value = 'Destination'
$ws.Cells.Item(11, 10).Value = 'Ops-linked flight'
$ws.Cells.Item(11, 11).Value = 'ET8'

# Row 12
$ws.Cells.Item(12, 1).Value = 'ZX-991122'
$ws.Cells.Item(12, 2).Value = 'Vehicle'
$ws.Cells.Item(12, 4).Value = 'Birmingham'
$ws.Cells.Item(12, 5).Value = 'B1 1AA'
$ws.Cells.Item(12, 8).Value = 'ACME Logistics Ltd'
$ws.Cells.Item(12, 9).Value = 'Company Vehicle'
$ws.Cells.Item(12, 10).Value = 'Fleet marker'
$ws.Cells.Item(12, 11).Value = 'ET3'

# Row 13
$ws.Cells.Item(13, 1).Value = 'AB12CDE'
$ws.Cells.Item(13, 2).Value = 'Vehicle'
$ws.Cells.Item(13, 4).Value = 'Leeds'
$ws.Cells.Item(13, 5).Value = 'LS1 4AP'
$ws.Cells.Item(13, 8).Value = 'Paul SHARP'
$ws.Cells.Item(13, 9).Value = 'Observed With'
$ws.Cells.Item(13, 10).Value = 'Vehicle observation'
$ws.Cells.Item(13, 11).Value = 'ET3'

# Row 14
$ws.Cells.Item(14, 1).Value = 'Moonshadow'
$ws.Cells.Item(14, 2).Value = 'Location'
$ws.Cells.Item(14, 4).Value = 'Marigot Bay Marina, Castries'
$ws.Cells.Item(14, 6).Value = 13.9697
$ws.Cells.Item(14, 7).Value = -61.0378
$ws.Cells.Item(14, 8).Value = 'Moonshadow Marine Ltd'
$ws.Cells.Item(14, 9).Value = 'Moored At'
$ws.Cells.Item(14, 10).Value = 'Vessel berth location'
$ws.Cells.Item(14, 11).Value = 'ET1'

# Row 15
$ws.Cells.Item(15, 1).Value = 'Waterfront Bar'
$ws.Cells.Item(15, 2).Value = 'Location'
$ws.Cells.Item(15, 4).Value = 'Harbour Front, Castries'
$ws.Cells.Item(15, 6).Value = 13.995
$ws.Cells.Item(15, 7).Value = -61.01
$ws.Cells.Item(15, 8).Value = 'Sam STEELE'
$ws.Cells.Item(15, 9).Value = 'Meeting'
$ws.Cells.Item(15, 10).Value = 'Witness meeting point'
$ws.Cells.Item(15, 11).Value = 'ET1'

# Row 16
$ws.Cells.Item(16, 1).Value = 'Sam STEELE'
$ws.Cells.Item(16, 2).Value = 'Person'
$ws.Cells.Item(16, 3).Value = "'" + '25/02/1984'
$ws.Cells.Item(16, 4).Value = 'Castries, Saint Lucia'
$ws.Cells.Item(16, 6).Value = 13.9972
$ws.Cells.Item(16, 7).Value = -61.0068
$ws.Cells.Item(16, 8).Value = 'Waterfront Bar'
$ws.Cells.Item(16, 9).Value = 'Sighted'
$ws.Cells.Item(16, 10).Value = 'Witness sighting'
$ws.Cells.Item(16, 11).Value = 'ET5'

# Row 17
$ws.Cells.Item(17, 1).Value = 'NCA-REF-2026-001'
$ws.Cells.Item(17, 2).Value = 'Communication'
$ws.Cells.Item(17, 8).Value = 'Andy JONES'
$ws.Cells.Item(17, 9).Value = 'Reference'
$ws.Cells.Item(17, 10).Value = 'Case reference handle'
$ws.Cells.Item(17, 11).Value = 'ET8'

# Row 18
$ws.Cells.Item(18, 1).Value = 'Finance Account 4931-1291'
$ws.Cells.Item(18, 2).Value = 'Organisation'
$ws.Cells.Item(18, 4).Value = 'Town Corp Bank, London'
$ws.Cells.Item(18, 5).Value = 'EC2M 7PP'
$ws.Cells.Item(18, 8).Value = 'ACME Logistics Ltd'
$ws.Cells.Item(18, 9).Value = 'Account Holder'
$ws.Cells.Item(18, 10).Value = 'Financial account proxy'
$ws.Cells.Item(18, 11).Value = 'ET4'

# Row 19
$ws.Cells.Item(19, 1).Value = 'EGDR - CULDROSE'
$ws.Cells.Item(19, 2).Value = 'Location'
$ws.Cells.Item(19, 4).Value = 'RNAS Culdrose, Helston'
$ws.Cells.Item(19, 5).Value = 'TR12 7RH'
$ws.Cells.Item(19, 8).Value = 'BAW23T'
$ws.Cells.Item(19, 9).Value = 'Destination'
$ws.Cells.Item(19, 10).Value = 'Airfield node'
$ws.Cells.Item(19, 11).Value = 'ET1'

# ---- Re-select N1 to match the saved view state ----
$ws.Range("N1").Select() | Out-Null

Write-Output "Import_Template expanded to 19 rows."
